# Update cryptos list (GitHub Actions data refresh) per commit:
# "Updated cryptos list on Tue Jun  4 04:46:13 UTC 2024 with GitHub Actions"
#
# Column D ("Price") cells are stored as plain text in the workbook (not
# numbers), even when their content happens to look numeric (e.g. "631.82").
# Assigning a numeric-looking string straight to Range.Value lets Excel's
# automatic type detection turn it into a real number, which would change
# the cell's stored type. To keep those cells as text - matching the
# source data - we briefly force a Text number format ("@") before writing
# the value, then restore the cell's style to "Normal" so no stray
# formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2: Bitcoin
$ws.Range("D2").Value = "69.207.54"
$ws.Range("E2").Value = "  +1.11%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.773.88"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.88%  "

# Row 5: BNB
Set-TextValue "D5" "631.82"
$ws.Range("E5").Value = "  +3.36%  "

# Row 6: Solana
Set-TextValue "D6" "166.92"
$ws.Range("E6").Value = "  +2.36%  "

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "3.771.49"
$ws.Range("E7").Value = "  -0.84%  "

# Row 8: USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9: XRP
$ws.Range("E9").Value = "  +0.93%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  -0.59%  "

# Row 11: Cardano
$ws.Range("E11").Value = "  +2.45%  "

# Row 12: Toncoin
Set-TextValue "D12" "6.76"
$ws.Range("E12").Value = "  -0.34%  "

# Row 13: ShibaInu
Set-TextValue "D13" "0.0000238"
$ws.Range("E13").Value = "  -3.34%  "

# Row 14: Avalanche
Set-TextValue "D14" "35.11"
$ws.Range("E14").Value = "  +0.37%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.409.78"
$ws.Range("E15").Value = "  -0.78%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "3.769.90"
$ws.Range("E16").Value = "  -2.04%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "69.217.04"
$ws.Range("E17").Value = "  +1.12%  "

# Row 18: Chainlink
Set-TextValue "D18" "17.62"
$ws.Range("E18").Value = "  -2.21%  "

# Row 19: TRON
$ws.Range("E19").Value = "  +0.14%  "

# Row 20: Polkadot
$ws.Range("E20").Value = "  -0.37%  "

# Row 21: BitcoinCash
Set-TextValue "D21" "462.75"
$ws.Range("E21").Value = "  +0.44%  "

# Row 22: Uniswap
$ws.Range("E22").Value = "  -0.54%  "

# Row 23: Polygon
$ws.Range("E23").Value = "  +1.46%  "

# Row 24: Litecoin
Set-TextValue "D24" "82.53"
$ws.Range("E24").Value = "  -0.85%  "

# Row 25: PEPE
$ws.Range("E25").Value = "  -1.28%  "

# Row 26: InternetComputer(DFINITY)
Set-TextValue "D26" "12.07"
$ws.Range("E26").Value = "  +1.26%  "

# Row 27: Fetch.AI
$ws.Range("E27").Value = "  +2.38%  "

# Row 28: RenderToken
Set-TextValue "D28" "10.08"
$ws.Range("E28").Value = "  +1.28%  "

# Row 29: Dai
$ws.Range("E29").Value = "  -0.08%  "

# Row 30: WrappedeETH
$ws.Range("D30").Value = "3.923.46"
$ws.Range("E30").Value = "  -0.68%  "

# Row 31: ImmutableX
Set-TextValue "D31" "2.32"
$ws.Range("E31").Value = "  +6.19%  "

# Row 32: PancakeSwap
$ws.Range("E32").Value = "  +2.65%  "

# Row 33: NEARProtocol
Set-TextValue "D33" "7.08"
$ws.Range("E33").Value = "  -1.36%  "

# Row 34 <-> Row 35: EthereumClassic and Kaspa swap ranking order, with fresh values
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D34" "0.179"
$ws.Range("E34").Value = "  +23.14%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D35" "28.48"
$ws.Range("E35").Value = "  -1.07%  "

# Row 36: Binance-PegBSC-USD
Set-TextValue "D36" "1.00"
$ws.Range("E36").Value = "  +0.00%  "

# Row 37: RenzoRestakedETH
$ws.Range("D37").Value = "3.726.33"
$ws.Range("E37").Value = "  -0.71%  "

# Row 38: Aptos
$ws.Range("E38").Value = "  -0.78%  "

# Row 39: Hedera
$ws.Range("E39").Value = "  +0.99%  "

# Row 40: dogwifhat
Set-TextValue "D40" "3.30"
$ws.Range("E40").Value = "  +5.71%  "

# Row 41: Filecoin
$ws.Range("E41").Value = "  -0.84%  "

# Row 42: FirstDigitalUSD
$ws.Range("E42").Value = "  -0.02%  "

# Row 43: Mantle
Set-TextValue "D43" "0.962"
$ws.Range("E43").Value = "  -1.66%  "

# Row 44: USDe
$ws.Range("E44").Value = "  -0.02%  "

# Row 45 <-> Row 46: Stacks and Monero swap ranking order, with fresh values
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D45" "158.18"
$ws.Range("E45").Value = "  +3.28%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D46" "1.99"
$ws.Range("E46").Value = "  +7.68%  "

# Row 47: ONDO
$ws.Range("E47").Value = "  +2.72%  "

# Row 48: Arweave
Set-TextValue "D48" "43.47"
$ws.Range("E48").Value = "  +1.39%  "

# Row 49 <-> Row 50: TheGraph and OKB swap ranking order, with fresh values
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D49" "47.11"
$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D50" "0.296"
$ws.Range("E50").Value = "  +0.39%  "

# Row 51: Cosmos
$ws.Range("E51").Value = "  +0.52%  "
